# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" columns for the
# 5538e567-8d54-4d41-b12e-f5a19a32c768 row (row 7) on both the
# "zh-cn" and "de-de" locale sheets, now that a handback has come in
# whose source commit is behind the latest source commit.

$wb = $excel.ActiveWorkbook

$sourceMd   = "5538e567-8d54-4d41-b12e-f5a19a32c768.md"
$sourceUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/614345cf1d226b297bef30b8c274f80cc037c1d5/e2e/5538e567-8d54-4d41-b12e-f5a19a32c768.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/decacb0c36a001d1c99d8c68f33c865cb6130f2a/e2e/5538e567-8d54-4d41-b12e-f5a19a32c768.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/614345cf1d226b297bef30b8c274f80cc037c1d5/e2e/5538e567-8d54-4d41-b12e-f5a19a32c768.md."

# ---- zh-cn sheet, row 7 ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value = "5538e567-8d54-4d41-b12e-f5a19a32c768.f6bd466b38ea355cf3e4263f0c8c1cad928c3a29.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-26 20:54:12"
$wsZh.Range("P7").Value = $errorDetail

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $sourceUrl, "", "", $sourceMd)
$wsZh.Range("I7").Style = "Hyperlink"

# ---- de-de sheet, row 7 ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value = "5538e567-8d54-4d41-b12e-f5a19a32c768.f6bd466b38ea355cf3e4263f0c8c1cad928c3a29.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-26 20:54:18"
$wsDe.Range("P7").Value = $errorDetail

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $sourceUrl, "", "", $sourceMd)
$wsDe.Range("I7").Style = "Hyperlink"
